$d = $word.ActiveDocument

# Step 1: Merge "Érvelések (Név, melyiket, miért" + "):" runs into a single run by
# replacing across the bookmark boundary. This also drops the stray _GoBack bookmark
# that currently sits mid-sentence (it will be re-added at the end of the new text).
$d.Content.Find.Execute("miért):", $true, $false, $false, $false, $false, $true, 1, $false, "miért):", 2) | Out-Null

# Step 2: Locate the "Érvelések..." paragraph and append the "Fecó, 2.-ik beteg:" line
# (preceded by a manual line break) to the same paragraph.
$argParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Érvelések")) {
        $argParaIndex = $i
        break
    }
}
$p = $d.Paragraphs.Item($argParaIndex)
$r = $p.Range
$r.MoveEnd(1, -1)
$r.Collapse(0)
$r.InsertBreak(6)
$r2 = $p.Range
$r2.MoveEnd(1, -1)
$r2.Collapse(0)
$r2.InsertAfter("Fecó, 2.-ik beteg:")

# Step 3: Insert a brand new paragraph after it with the argument text for Fecó.
$p = $d.Paragraphs.Item($argParaIndex)
$pRange = $p.Range
$pRange.Collapse(0)
$pRange.InsertParagraphAfter()

$newP = $d.Paragraphs.Item($argParaIndex + 1)
$nr = $newP.Range
$nr.MoveEnd(1, -1)
$nr.Collapse(0)
$nr.InsertAfter("`tAmennyiben sikerül a műtét a polgármester képes lenne több pénzel támogatni a Kórházat ezzel több beteget megmentve. A 4-ik semmi képpen sem lenne jó mivel csak kihasználná a kórházat, és nem tesz hozzá semmit a pénzével. Az első habár ki bírná fizetni a műtétet anyagi helyzetét nem lenne képes fenttartani. A 3.-ik habár képes lenne anyagilag tovább élni nem tudna nagy változást tenni a kórház anyagi helyzetén.")

Write-Host "Done"
